$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

# Cells A3/B3/C3 use a cell style with quotePrefix="1" set. A plain
# .Value assignment re-evaluates the cell and drops that quote-prefix
# flag (the style index would silently shift to the sibling style that
# is identical except for quotePrefix="0"). To keep the original
# formatting intact we stash the existing format in an unused scratch
# cell, write the new values, then paste the stashed format back.
$scratch = $ws.Range("ZZ100")
$ws.Range("A3").Copy($scratch)

$ws.Range("A3").Value = "JSSO1000243"
$ws.Range("B3").Value = "JSSO1000243"
$ws.Range("C3").Value = "JSSO1000243"

$scratch.Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$scratch.Clear()

# Remaining changed cells keep the same cell style before and after
# the edit, so a direct value assignment is sufficient.
$ws.Range("AJ3").Value = "JSCN1000243"
$ws.Range("AL3").Value = "SLJSSO1000243"
$ws.Range("AN3").Value = "JSCN1000243"
